$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Helper: locate the 1-based index of the paragraph whose range
# text equals an exact string (after trimming the trailing
# paragraph-mark / cell-mark characters Word includes in
# Range.Text). Re-run this fresh after every mutation -- Paragraph
# / Range object references captured before an insert/delete can
# end up "stuck" at their original slot, so we always re-resolve
# by index.
# ---------------------------------------------------------------
function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# =================================================================
# CHANGE 4 (performed first -- it is the last edit location in the
# document, so handling it first keeps every earlier paragraph
# index stable for the edits that follow).
# Insert 3 new paragraphs between "{Image_1}" and "{Image_2}":
#   empty / "{Table_2}" / empty
# =================================================================
$idxImage1 = Get-ParaIndexByText $d "{Image_1}"
$d.Paragraphs.Item($idxImage1).Range.InsertParagraphAfter()

$idxEmptyA = $idxImage1 + 1
$d.Paragraphs.Item($idxEmptyA).Range.InsertParagraphAfter()

$idxTable2 = $idxEmptyA + 1
$d.Paragraphs.Item($idxTable2).Range.Text = "{Table_2}"

$d.Paragraphs.Item($idxTable2).Range.InsertParagraphAfter()

Write-Output "Change 4 done"

# =================================================================
# CHANGE 3
# Of the 4 empty paragraphs right before "{Graph_Eq}":
#   - keep the first one
#   - delete the 2nd and 3rd
#   - insert a new paragraph "Obtained Linear Fite Equation" (bold+underline)
#   - delete the 4th, moving its "center" alignment onto "{Graph_Eq}"
# =================================================================
$idxGraphEq = Get-ParaIndexByText $d "{Graph_Eq}"
$idxEmpty4 = $idxGraphEq - 1
$idxEmpty3 = $idxGraphEq - 2
$idxEmpty2 = $idxGraphEq - 3
$idxEmpty1 = $idxGraphEq - 4

# delete 3rd then 2nd (descending order keeps earlier indices valid)
$d.Paragraphs.Item($idxEmpty3).Range.Delete()
$d.Paragraphs.Item($idxEmpty2).Range.Delete()

# after those two deletions, the 1st empty paragraph is still at $idxEmpty1
# and the 4th empty paragraph / {Graph_Eq} shifted down by 2.
$d.Paragraphs.Item($idxEmpty1).Range.InsertParagraphAfter()
$idxNewTitle = $idxEmpty1 + 1
$pNewTitle = $d.Paragraphs.Item($idxNewTitle)
$pNewTitle.Range.Text = "Obtained Linear Fite Equation"
$pNewTitle.Range.Font.Bold = 1
$pNewTitle.Range.Font.BoldBi = 1
$pNewTitle.Range.Font.Underline = 1

# re-resolve fresh indices (an extra paragraph was just inserted)
$idxGraphEq2 = Get-ParaIndexByText $d "{Graph_Eq}"
$idxEmpty4b = $idxGraphEq2 - 1
$d.Paragraphs.Item($idxGraphEq2).Alignment = 1
$d.Paragraphs.Item($idxEmpty4b).Range.Delete()

Write-Output "Change 3 done"

# =================================================================
# CHANGE 2
# "{Table}" -> "{Table_1}"
# =================================================================
$d.Content.Find.Execute("{Table}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{Table_1}", 2) | Out-Null

Write-Output "Change 2 done"

# =================================================================
# CHANGE 1
# Insert a new paragraph "Normalised Equation" (bold + underline)
# right after the "Bradford Assay" heading paragraph.
# =================================================================
$idxHeading = Get-ParaIndexByText $d "Bradford Assay"
$idxNext = $idxHeading + 1
$d.Paragraphs.Item($idxNext).Range.InsertParagraphBefore()

$pTitle = $d.Paragraphs.Item($idxNext)
$pTitle.Range.Text = "Normalised Equation"
$pTitle.Range.Font.Bold = 1
$pTitle.Range.Font.BoldBi = 1
$pTitle.Range.Font.Underline = 1

Write-Output "Change 1 done"

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
